# Apply the "Diagramme avec References" commit:
#  - Update the Metadata "Date" value (sheet "Metadata", B8)
#  - Add two new FHIR logical-model elements ("ExerciceProfessionnel" and
#    "EntiteGeographique" reference links) as rows 16 and 17 on the
#    "Elements" sheet, following the same shape as the existing rows.
#  - Widen a few columns to fit the new, slightly longer, content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the generation Date stamp.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: append the two new element rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Give the two new rows the same formatting as the last existing data row
# (row 15) before filling in values.
$ws.Range("A16:AJ16").Style = $ws.Range("A15:AJ15").Style
$ws.Range("A17:AJ17").Style = $ws.Range("A15:AJ15").Style

# Note: the "Min"/"Max"/"Base Min"/"Base Max" columns hold textual "0"/"1"
# labels (not numbers) in this table, same as the existing rows above.
# NumberFormat="@" forces text storage for the otherwise-numeric-looking
# literal; re-applying the row-15 cell style afterwards restores the
# normal look (NumberFormat="@" would otherwise leave behind a one-off
# style that drops the cell's border/fill).

# --- Row 16: Certificat.ExerciceProfessionnel -----------------------------
$ws.Cells.Item(16, 1).Value = "Certificat.ExerciceProfessionnel"
$ws.Cells.Item(16, 2).Value = "Certificat.ExerciceProfessionnel"
$ws.Cells.Item(16, 4).Value = ""

$ws.Cells.Item(16, 6).NumberFormat = "@"
$ws.Cells.Item(16, 6).Value = "0"
$ws.Cells.Item(16, 6).Style = $ws.Cells.Item(15, 6).Style

$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "1"
$ws.Cells.Item(16, 7).Style = $ws.Cells.Item(15, 7).Style

$ws.Cells.Item(16, 8).Value = ""
$ws.Cells.Item(16, 9).Value = ""
$ws.Cells.Item(16, 10).Value = ""
$ws.Cells.Item(16, 11).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/ExerciceProfessionnel`n"
$ws.Cells.Item(16, 12).Value = "Lien vers la classe ExerciceProfessionnel"
$ws.Cells.Item(16, 13).Value = "Lien vers la classe ExerciceProfessionnel"
$ws.Cells.Item(16, 16).Value = ""
$ws.Cells.Item(16, 18).Value = ""
$ws.Cells.Item(16, 19).Value = ""
$ws.Cells.Item(16, 20).Value = ""
$ws.Cells.Item(16, 21).Value = ""
$ws.Cells.Item(16, 22).Value = ""
$ws.Cells.Item(16, 23).Value = ""
$ws.Cells.Item(16, 24).Value = ""
$ws.Cells.Item(16, 25).Value = ""
$ws.Cells.Item(16, 26).Value = ""
$ws.Cells.Item(16, 27).Value = ""
$ws.Cells.Item(16, 28).Value = ""
$ws.Cells.Item(16, 29).Value = ""
$ws.Cells.Item(16, 30).Value = ""
$ws.Cells.Item(16, 31).Value = ""
$ws.Cells.Item(16, 32).Value = "Certificat.ExerciceProfessionnel"

$ws.Cells.Item(16, 33).NumberFormat = "@"
$ws.Cells.Item(16, 33).Value = "0"
$ws.Cells.Item(16, 33).Style = $ws.Cells.Item(15, 33).Style

$ws.Cells.Item(16, 34).NumberFormat = "@"
$ws.Cells.Item(16, 34).Value = "1"
$ws.Cells.Item(16, 34).Style = $ws.Cells.Item(15, 34).Style

$ws.Cells.Item(16, 35).Value = ""
$ws.Cells.Item(16, 36).Value = ""

# --- Row 17: Certificat.EntiteGeographique --------------------------------
$ws.Cells.Item(17, 1).Value = "Certificat.EntiteGeographique"
$ws.Cells.Item(17, 2).Value = "Certificat.EntiteGeographique"
$ws.Cells.Item(17, 4).Value = ""

$ws.Cells.Item(17, 6).NumberFormat = "@"
$ws.Cells.Item(17, 6).Value = "0"
$ws.Cells.Item(17, 6).Style = $ws.Cells.Item(15, 6).Style

$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "1"
$ws.Cells.Item(17, 7).Style = $ws.Cells.Item(15, 7).Style

$ws.Cells.Item(17, 8).Value = ""
$ws.Cells.Item(17, 9).Value = ""
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 11).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteGeographique`n"
$ws.Cells.Item(17, 12).Value = "Lien vers la classe EntiteGeographique"
$ws.Cells.Item(17, 13).Value = "Lien vers la classe EntiteGeographique"
$ws.Cells.Item(17, 16).Value = ""
$ws.Cells.Item(17, 18).Value = ""
$ws.Cells.Item(17, 19).Value = ""
$ws.Cells.Item(17, 20).Value = ""
$ws.Cells.Item(17, 21).Value = ""
$ws.Cells.Item(17, 22).Value = ""
$ws.Cells.Item(17, 23).Value = ""
$ws.Cells.Item(17, 24).Value = ""
$ws.Cells.Item(17, 25).Value = ""
$ws.Cells.Item(17, 26).Value = ""
$ws.Cells.Item(17, 27).Value = ""
$ws.Cells.Item(17, 28).Value = ""
$ws.Cells.Item(17, 29).Value = ""
$ws.Cells.Item(17, 30).Value = ""
$ws.Cells.Item(17, 31).Value = ""
$ws.Cells.Item(17, 32).Value = "Certificat.EntiteGeographique"

$ws.Cells.Item(17, 33).NumberFormat = "@"
$ws.Cells.Item(17, 33).Value = "0"
$ws.Cells.Item(17, 33).Style = $ws.Cells.Item(15, 33).Style

$ws.Cells.Item(17, 34).NumberFormat = "@"
$ws.Cells.Item(17, 34).Value = "1"
$ws.Cells.Item(17, 34).Style = $ws.Cells.Item(15, 34).Style

$ws.Cells.Item(17, 35).Value = ""
$ws.Cells.Item(17, 36).Value = ""

# ---------------------------------------------------------------------------
# 3. Widen the columns whose best-fit width grew because of the new,
#    slightly longer cell contents (ID/Path/Base Path + Type(s) columns).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.0
$ws.Columns.Item(2).ColumnWidth = 25.0
$ws.Columns.Item(11).ColumnWidth = 62.166666666666664
$ws.Columns.Item(32).ColumnWidth = 25.0
